# Auto-generated edit script: Add data for 2024-07-17
# Applies updated 2024 (column K, and a couple column J) violent crime counts
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 4296
$ws.Range("K3").Value = 4363
$ws.Range("K5").Value = 321
$ws.Range("K6").Value = 4897
$ws.Range("K7").Value = 14758

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 129
$ws.Range("K8").Value = 992
$ws.Range("K11").Value = 293
$ws.Range("K18").Value = 101
$ws.Range("K19").Value = 447
$ws.Range("K20").Value = 335
$ws.Range("K22").Value = 43
$ws.Range("K24").Value = 44
$ws.Range("K25").Value = 69
$ws.Range("K29").Value = 784
$ws.Range("K30").Value = 53
$ws.Range("K31").Value = 161
$ws.Range("K33").Value = 617
$ws.Range("K37").Value = 505
$ws.Range("K42").Value = 543
$ws.Range("K43").Value = 129
$ws.Range("K47").Value = 91
$ws.Range("K48").Value = 191
$ws.Range("K52").Value = 399
$ws.Range("K53").Value = 198
$ws.Range("J54").Value = 574
$ws.Range("K54").Value = 274
$ws.Range("K55").Value = 166
$ws.Range("K60").Value = 96
$ws.Range("J63").Value = 108
$ws.Range("K63").Value = 42
$ws.Range("K67").Value = 569
$ws.Range("K68").Value = 38
$ws.Range("K70").Value = 24
$ws.Range("K72").Value = 67
$ws.Range("K78").Value = 175
$ws.Range("K79").Value = 378
$ws.Range("K83").Value = 313
$ws.Range("K84").Value = 107
$ws.Range("K85").Value = 664
$ws.Range("K86").Value = 100
$ws.Range("K88").Value = 169
$ws.Range("K89").Value = 205
$ws.Range("K90").Value = 132
$ws.Range("K94").Value = 184
$ws.Range("K95").Value = 259
$ws.Range("K96").Value = 162
$ws.Range("K97").Value = 124
$ws.Range("K101").Value = 14758

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K6").Value = 73
$ws.Range("K7").Value = 162

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 95
$ws.Range("K3").Value = 74
$ws.Range("K6").Value = 105
$ws.Range("K7").Value = 293

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K2").Value = 56
$ws.Range("K6").Value = 67
$ws.Range("K7").Value = 205

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 235
$ws.Range("K3").Value = 220
$ws.Range("K4").Value = 35
$ws.Range("K5").Value = 19
$ws.Range("K7").Value = 664

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 104
$ws.Range("K7").Value = 399

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value = 49
$ws.Range("K7").Value = 198

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 278
$ws.Range("K3").Value = 296
$ws.Range("K6").Value = 334
$ws.Range("K7").Value = 992

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K3").Value = 109
$ws.Range("K7").Value = 313

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 169
$ws.Range("K3").Value = 230
$ws.Range("K4").Value = 26
$ws.Range("K6").Value = 178
$ws.Range("K7").Value = 617

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 85
$ws.Range("K3").Value = 89
$ws.Range("K6").Value = 62
$ws.Range("K7").Value = 259

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K2").Value = 139
$ws.Range("K3").Value = 165
$ws.Range("K6").Value = 152
$ws.Range("K7").Value = 505

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 57
$ws.Range("K7").Value = 161

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K6").Value = 167
$ws.Range("K7").Value = 569

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 107

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 76
$ws.Range("J4").Value = 49
$ws.Range("K4").Value = 14
$ws.Range("J7").Value = 574
$ws.Range("K7").Value = 274

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 225
$ws.Range("K3").Value = 278
$ws.Range("K6").Value = 217
$ws.Range("K7").Value = 784

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K3").Value = 45
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 135
$ws.Range("K6").Value = 137
$ws.Range("K7").Value = 447

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 150
$ws.Range("K6").Value = 197
$ws.Range("K7").Value = 543

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 51
$ws.Range("K6").Value = 64
$ws.Range("K7").Value = 175

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K2").Value = 48
$ws.Range("K3").Value = 46
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 166

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 126
$ws.Range("K3").Value = 123
$ws.Range("K7").Value = 378

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 117
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 335

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K3").Value = 52
$ws.Range("K6").Value = 43

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K2").Value = 54
$ws.Range("K3").Value = 31
$ws.Range("K7").Value = 184

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 69

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K3").Value = 27
$ws.Range("K4").Value = 6
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 91

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 39
$ws.Range("K5").Value = 6
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("K3").Value = 23
$ws.Range("K7").Value = 124

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 53
$ws.Range("K7").Value = 169

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K4").Value = 39
$ws.Range("K7").Value = 100

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 46
$ws.Range("K6").Value = 27
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("K2").Value = 16
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 96

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 129

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K2").Value = 22
$ws.Range("K7").Value = 43

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 67
